# Apply marksheet corrections: update correct/total marks figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row (B11): corrected marks value 3 -> 5
$ws.Range("B11").Value = 5

# Total row (B12): corrected marks value 54 -> 90
$ws.Range("B12").Value = 90

# Total row (E12): corrected fraction text "47/84" -> "90/140"
$ws.Range("E12").Value = "90/140"
